# Apply change described in commit:
# "단위값에 맞추어 임시값 작성 Character - MovementSpd 5.75 -> 575"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CharacterGameData")

# Update MovementSpd value (cell F3) from 5.75 to 575
$ws.Range("F3").Value = 575

# Move the active selection to I6 (matches the selection state captured in the saved file)
$ws.Activate()
$ws.Range("I6").Select()
